# competition.xlsx update — "updated gps, random particles, and starting to
# integrate compass"
#
# This adds a handful of new trajectory way-points (with a new "wall
# following" flag column) to the traj sheet, shifts the final M2 marker row
# down to make room for them, and updates the map sheet's selection to match
# the author's last cursor position.

$wb   = $excel.ActiveWorkbook
$map  = $wb.Worksheets.Item("map")
$traj = $wb.Worksheets.Item("traj")

# --- traj sheet: new "wall following" column header ---------------------
$traj.Range("E1").Value = "wall following"

# --- traj sheet: new way-points (rows 20, 22, 24) ------------------------
$traj.Range("B20").Value = 20
$traj.Range("C20").Value = -26

$traj.Range("B22").Value = 20
$traj.Range("C22").Value = -21
$traj.Range("E22").Value = 1

$traj.Range("B24").Value = 20
$traj.Range("C24").Value = -9

# --- traj sheet: row 26 keeps the trajectory going (was the old marker
#     row), and the "M2" end marker now lives on row 28 -------------------
$traj.Range("A26").ClearContents() | Out-Null
$traj.Range("B26").Value = 50
$traj.Range("C26").Value = -9

$traj.Range("A28").Value = "M2"
$traj.Range("B28").Value = 50
$traj.Range("C28").Value = -27

# --- map sheet: cursor moved to D19 ---------------------------------------
$map.Range("D19").Select() | Out-Null
